# Auto-generated Excel COM-interop script
# Updates market-price derived columns (H:N) on rows across the 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed market data
# pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3018.5715
$ws.Range("I19").Value = 6174.778
$ws.Range("J19").Value = 651.4167
$ws.Range("K19").Value = 6174.778
$ws.Range("L19").Value = 651.4167
$ws.Range("M19").Value = -5999.778
$ws.Range("N19").Value = -1001.4167
$ws.Range("H113").Value = 2791.7104
$ws.Range("I113").Value = 2342.5881
$ws.Range("J113").Value = 3155.2856
$ws.Range("K113").Value = 2342.5881
$ws.Range("L113").Value = 3155.2856
$ws.Range("M113").Value = 911.4119000000001
$ws.Range("N113").Value = -9663.285599999999
$ws.Range("H132").Value = 3450.2632
$ws.Range("I132").Value = 1824.6595
$ws.Range("J132").Value = 11090.6
$ws.Range("K132").Value = 5473.9785
$ws.Range("L132").Value = 33271.8
$ws.Range("M132").Value = -2943.9785
$ws.Range("N132").Value = -38331.8
$ws.Range("H137").Value = 3485.5334
$ws.Range("I137").Value = 3752.85
$ws.Range("J137").Value = 2950.9
$ws.Range("K137").Value = 11258.55
$ws.Range("L137").Value = 8852.700000000001
$ws.Range("M137").Value = -8708.549999999999
$ws.Range("N137").Value = -13952.7
$ws.Range("H141").Value = 2050.8823
$ws.Range("I141").Value = 2050.8823
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6152.646900000001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -972.6469000000006
$ws.Range("N141").ClearContents()

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2897.7778
$ws.Range("I2").Value = 4540
$ws.Range("J2").Value = 845
$ws.Range("K2").Value = 4540
$ws.Range("L2").Value = 845
$ws.Range("M2").Value = -4427
$ws.Range("N2").Value = -1071
$ws.Range("H45").Value = 1476.125
$ws.Range("I45").Value = 1190
$ws.Range("J45").Value = 1762.25
$ws.Range("K45").Value = 1190
$ws.Range("L45").Value = 1762.25
$ws.Range("M45").Value = -813
$ws.Range("N45").Value = -2516.25
$ws.Range("H110").Value = 1511.9131
$ws.Range("I110").Value = 1389.7273
$ws.Range("K110").Value = 1389.7273
$ws.Range("M110").Value = 655.2727
$ws.Range("H116").Value = 2897.7778
$ws.Range("I116").Value = 4540
$ws.Range("J116").Value = 845
$ws.Range("K116").Value = 4540
$ws.Range("L116").Value = 845
$ws.Range("M116").Value = -2246
$ws.Range("N116").Value = -5433
$ws.Range("H141").Value = 29582.25
$ws.Range("J141").Value = 29582.25
$ws.Range("L141").Value = 29582.25
$ws.Range("N141").Value = -39942.25

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2897.7778
$ws.Range("I3").Value = 4540
$ws.Range("J3").Value = 845
$ws.Range("K3").Value = 4540
$ws.Range("L3").Value = 845
$ws.Range("M3").Value = -4426
$ws.Range("N3").Value = -1073
$ws.Range("H134").Value = 17350.148
$ws.Range("I134").Value = 21705.428
$ws.Range("J134").Value = 5494.1113
$ws.Range("K134").Value = 65116.284
$ws.Range("L134").Value = 16482.3339
$ws.Range("M134").Value = -62581.284
$ws.Range("N134").Value = -21552.3339

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2782.6316
$ws.Range("I31").Value = 2060.6365
$ws.Range("J31").Value = 3775.375
$ws.Range("K31").Value = 2060.6365
$ws.Range("L31").Value = 3775.375
$ws.Range("M31").Value = -1765.6365
$ws.Range("N31").Value = -4365.375
$ws.Range("H34").Value = 2782.6316
$ws.Range("I34").Value = 2060.6365
$ws.Range("J34").Value = 3775.375
$ws.Range("K34").Value = 2060.6365
$ws.Range("L34").Value = 3775.375
$ws.Range("M34").Value = -1858.6365
$ws.Range("N34").Value = -4179.375
$ws.Range("H107").Value = 527
$ws.Range("I107").Value = 415.7143
$ws.Range("J107").Value = 682.8
$ws.Range("K107").Value = 415.7143
$ws.Range("L107").Value = 682.8
$ws.Range("M107").Value = 1504.2857
$ws.Range("N107").Value = -4522.8
$ws.Range("H122").Value = 1197.6
$ws.Range("I122").Value = 870.8570999999999
$ws.Range("J122").Value = 1960
$ws.Range("K122").Value = 2612.5713
$ws.Range("L122").Value = 5880
$ws.Range("M122").Value = -162.5712999999996
$ws.Range("N122").Value = -10780
$ws.Range("H132").Value = 2190.717
$ws.Range("I132").Value = 1364.6177
$ws.Range("J132").Value = 3669
$ws.Range("K132").Value = 4093.8531
$ws.Range("L132").Value = 11007
$ws.Range("M132").Value = -1563.8531
$ws.Range("N132").Value = -16067

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2554.12
$ws.Range("J69").Value = 3550
$ws.Range("L69").Value = 10650
$ws.Range("N69").Value = -12272
$ws.Range("H72").Value = 2554.12
$ws.Range("J72").Value = 3550
$ws.Range("L72").Value = 31950
$ws.Range("N72").Value = -40062
$ws.Range("H132").Value = 4728.5713
$ws.Range("I132").Value = 2225
$ws.Range("J132").Value = 5317.647
$ws.Range("K132").Value = 20025
$ws.Range("L132").Value = 47858.823
$ws.Range("M132").Value = -17495
$ws.Range("N132").Value = -52918.823

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7699.606
$ws.Range("I70").Value = 3876.4614
$ws.Range("J70").Value = 21899.857
$ws.Range("K70").Value = 3876.4614
$ws.Range("L70").Value = 21899.857
$ws.Range("M70").Value = -3606.4614
$ws.Range("N70").Value = -22439.857
$ws.Range("H73").Value = 7699.606
$ws.Range("I73").Value = 3876.4614
$ws.Range("J73").Value = 21899.857
$ws.Range("K73").Value = 3876.4614
$ws.Range("L73").Value = 21899.857
$ws.Range("M73").Value = -2940.4614
$ws.Range("N73").Value = -23771.857
$ws.Range("H80").Value = 4402.931
$ws.Range("I80").Value = 5222.5
$ws.Range("J80").Value = 3061.818
$ws.Range("K80").Value = 5222.5
$ws.Range("L80").Value = 3061.818
$ws.Range("M80").Value = -4224.5
$ws.Range("N80").Value = -5057.818
$ws.Range("H83").Value = 4402.931
$ws.Range("I83").Value = 5222.5
$ws.Range("J83").Value = 3061.818
$ws.Range("K83").Value = 26112.5
$ws.Range("L83").Value = 15309.09
$ws.Range("M83").Value = -21120.5
$ws.Range("N83").Value = -25293.09
$ws.Range("H102").Value = 3023.8667
$ws.Range("I102").Value = 3059.1724
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 3059.1724
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -1437.1724
$ws.Range("N102").Value = -5244
$ws.Range("H113").Value = 11166.6
$ws.Range("I113").Value = 14960.286
$ws.Range("J113").Value = 2314.6667
$ws.Range("K113").Value = 14960.286
$ws.Range("L113").Value = 2314.6667
$ws.Range("M113").Value = -12790.286
$ws.Range("N113").Value = -6654.6667
$ws.Range("H132").Value = 3709.2173
$ws.Range("I132").Value = 3796.5925
$ws.Range("J132").Value = 3585.0527
$ws.Range("K132").Value = 11389.7775
$ws.Range("L132").Value = 10755.1581
$ws.Range("M132").Value = -8859.7775
$ws.Range("N132").Value = -15815.1581

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 926.6667
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 990
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 990
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -1366
$ws.Range("H132").Value = 5794.375
$ws.Range("I132").Value = 2141.6875
$ws.Range("J132").Value = 10664.625
$ws.Range("K132").Value = 6425.0625
$ws.Range("L132").Value = 31993.875
$ws.Range("M132").Value = -3895.0625
$ws.Range("N132").Value = -37053.875

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 557.5
$ws.Range("I107").Value = 555.0909
$ws.Range("J107").Value = 560.44446
$ws.Range("K107").Value = 1665.2727
$ws.Range("L107").Value = 1681.33338
$ws.Range("M107").Value = 254.7273
$ws.Range("N107").Value = -5521.33338
$ws.Range("H122").Value = 37620.82
$ws.Range("I122").Value = 45246.87
$ws.Range("K122").Value = 135740.61
$ws.Range("M122").Value = -133290.61
$ws.Range("H136").Value = 1824.5
$ws.Range("I136").Value = 1072.8529
$ws.Range("J136").Value = 2807.423
$ws.Range("K136").Value = 3218.5587
$ws.Range("L136").Value = 8422.269
$ws.Range("M136").Value = -668.5587000000005
$ws.Range("N136").Value = -13522.269

Write-Output "Updated market-price columns across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
